$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New trade rows added in the 1-Jun-2021 midday update.
# Each pair of rows (SELL leg + BUY leg) represents one option spread trade;
# a single blank row separates consecutive trades (matching existing layout).
$newRows = @(
    @{ Row=475; A=44340; B="1511-87668"; C="SELL"; D=100; E="SPX"; F=44372; G="PUT"; H=4090; I=43.92 },
    @{ Row=476; A=44340; B="1511-87668"; C="BUY"; D=-100; E="SPX"; F=44372; G="PUT"; H=4085; I=43.07 },
    @{ Row=478; A=44340; B="1511-88542"; C="SELL"; D=100; E="SPX"; F=44354; G="PUT"; H=4090; I=16.83 },
    @{ Row=479; A=44340; B="1511-88542"; C="BUY"; D=-100; E="SPX"; F=44354; G="PUT"; H=4085; I=16.28 },
    @{ Row=481; A=44340; B="1511-96435"; C="BUY"; D=-100; E="RUT"; F=44372; G="PUT"; H=2115; I=27.8 },
    @{ Row=482; A=44340; B="1511-96435"; C="SELL"; D=100; E="RUT"; F=44372; G="PUT"; H=2120; I=28.84 },
    @{ Row=484; A=44341; B="1514-10989"; C="BUY"; D=-100; E="PUT"; F=44372; G="PUT"; H=2120; I=23.49 },
    @{ Row=485; A=44341; B="1514-10989"; C="SELL"; D=100; E="PUT"; F=44372; G="PUT"; H=2105; I=20.82 },
    @{ Row=487; A=44341; B="1514-17535"; C="BUY"; D=-100; E="RUT"; F=44372; G="CALL"; H=2305; I=27.2 },
    @{ Row=488; A=44341; B="1514-17535"; C="SELL"; D=100; E="RUT"; F=44372; G="CALL"; H=2310; I=25.4 },
    @{ Row=490; A=44341; B="1514-19827"; C="SELL"; D=100; E="SPX"; F=44354; G="PUT"; H=4100; I=14.62 },
    @{ Row=491; A=44341; B="1514-19827"; C="BUY"; D=-100; E="SPX"; F=44354; G="PUT"; H=4090; I=13.47 },
    @{ Row=493; A=44341; B="1514-24734"; C="SELL"; D=100; E="SPX"; F=44372; G="PUT"; H=4095; I=42.47 },
    @{ Row=494; A=44341; B="1514-24734"; C="BUY"; D=-100; E="SPX"; F=44372; G="PUT"; H=4090; I=41.57 },
    @{ Row=496; A=44342; B="1517-12588"; C="SELL"; D=100; E="SPX"; F=44393; G="PUT"; H=4000; I=46.22 },
    @{ Row=497; A=44342; B="1517-12588"; C="BUY"; D=-100; E="SPX"; F=44393; G="PUT"; H=3980; I=43.32 },
    @{ Row=499; A=44342; B="1517-18557"; C="SELL"; D=100; E="NDX"; F=44379; G="PUT"; H=13100; I=170.75 },
    @{ Row=500; A=44342; B="1517-18557"; C="BUY"; D=-100; E="NDX"; F=44379; G="PUT"; H=13075; I=165.95 },
    @{ Row=502; A=44343; B="1519-35625"; C="SELL"; D=100; E="SPX"; F=44372; G="PUT"; H=4100; I=34.87 },
    @{ Row=503; A=44343; B="1519-35625"; C="BUY"; D=-100; E="SPX"; F=44372; G="PUT"; H=4095; I=34.07 },
    @{ Row=505; A=44343; B="1519-39131"; C="BUY"; D=-100; E="SPX"; F=44354; G="PUT"; H=4100; I=9.08 },
    @{ Row=506; A=44343; B="1519-39131"; C="SELL"; D=100; E="SPX"; F=44354; G="PUT"; H=4095; I=8.63 },
    @{ Row=508; A=44343; B="1519-40552"; C="BUY"; D=-100; E="SPX"; F=44393; G="PUT"; H=4000; I=39.63 },
    @{ Row=509; A=44343; B="1519-40552"; C="SELL"; D=100; E="SPX"; F=44393; G="PUT"; H=3990; I=38.23 },
    @{ Row=511; A=44343; B="1519-65314"; C="BUY"; D=-100; E="RUT"; F=44372; G="PUT"; H=2105; I=15.18 },
    @{ Row=512; A=44343; B="1519-65314"; C="SELL"; D=100; E="RUT"; F=44372; G="PUT"; H=2110; I=15.8 },
    @{ Row=514; A=44343; B="1519-65891"; C="BUY"; D=-100; E="RUT"; F=44372; G="CALL"; H=2310; I=29.92 },
    @{ Row=515; A=44343; B="1519-65891"; C="SELL"; D=100; E="RUT"; F=44372; G="CALL"; H=2305; I=31.97 },
    @{ Row=517; A=44344; B="1522-52805"; C="SELL"; D=100; E="SPX"; F=44393; G="PUT"; H=4070; I=45.32 },
    @{ Row=518; A=44344; B="1522-52805"; C="BUY"; D=-100; E="SPX"; F=44393; G="PUT"; H=3990; I=33.82 },
    @{ Row=520; A=44344; B="1522-66478"; C="SELL"; D=100; E="SPX"; F=44354; G="PUT"; H=4150; I=10.95 },
    @{ Row=521; A=44344; B="1522-66478"; C="BUY"; D=-100; E="SPX"; F=44354; G="PUT"; H=4095; I=5.85 },
    @{ Row=523; A=44344; B="1522-94554"; C="SELL"; D=100; E="RUT"; F=44372; G="PUT"; H=2100; I=13.02 },
    @{ Row=524; A=44344; B="1522-94554"; C="BUY"; D=-100; E="RUT"; F=44372; G="PUT"; H=2070; I=10.17 },
    @{ Row=526; A=44344; B="1522-98714"; C="SELL"; D=100; E="RUT"; F=44372; G="CALL"; H=2310; I=29.78 },
    @{ Row=527; A=44344; B="1522-98714"; C="BUY"; D=-100; E="RUT"; F=44372; G="CALL"; H=2340; I=19.03 },
    @{ Row=529; A=44344; B="1523-02944"; C="SELL"; D=100; E="SPX"; F=44372; G="PUT"; H=4075; I=26.69 },
    @{ Row=530; A=44344; B="1523-02944"; C="BUY"; D=-100; E="SPX"; F=44372; G="PUT"; H=4020; I=20.44 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
}

# Write the TOTAL column (K) formulas. Each trade occupies two consecutive rows,
# so filling the formula across each two-row block reproduces Excel's shared-formula
# grouping (same mechanism already used for the existing K445:K473 block).
for ($i = 0; $i -lt $newRows.Count; $i += 2) {
    $r1 = $newRows[$i].Row
    $r2 = $newRows[$i + 1].Row
    $ws.Range("K$($r1):K$($r2)").Formula = "=D$r1*I$r1"
}

# Move the selection to the new bottom of the sheet, matching the post-edit state.
$ws.Range("K532").Select()
